$ws = $excel.ActiveSheet

# Force text formatting on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (stripping meaningful trailing zeros, etc.)
$textCells = @("D5","D6","D7","D9","D10","D11","D13","D14","D15","D17","D20","D22","D23","D25","D27","D28","D30","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin / link / price / volume values
$ws.Range("D2").Value = "36.310.08"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.986.76"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "245.21"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").Value = "62.67"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("D10").Value = "56.42"
$ws.Range("E10").Value = "  -4.19%  "
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  +8.97%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "0.871"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").Value = "22.60"
$ws.Range("E14").Value = "  +11.46%  "
$ws.Range("D15").Value = "14.10"
$ws.Range("E15").Value = "  -5.70%  "
$ws.Range("D16").Value = "2.277.89"
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("D17").Value = "5.47"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "1.994.85"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "36.178.67"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("D20").Value = "71.39"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "0.0₃0876"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "5.29"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "237.08"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  -10.48%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "9.89"
$ws.Range("E27").Value = "  +4.12%  "
$ws.Range("D28").Value = "160.07"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  +22.53%  "
$ws.Range("D30").Value = "19.97"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "4.94"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").Value = "1.15"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("D34").Value = "0.0627"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").Value = "4.42"
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").Value = "6.35"
$ws.Range("E36").Value = "  +6.31%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +14.97%  "
$ws.Range("D41").Value = "0.0996"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").Value = "1.24"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "2.85"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "93.66"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "16.40"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "7.59"
$ws.Range("E48").Value = "  -5.85%  "
$ws.Range("D49").Value = "1.356.49"
$ws.Range("E49").Value = "  -4.70%  "
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "2.174.32"
$ws.Range("E51").Value = "  -1.80%  "
